# Add a new log row (row 4) documenting the 2023-08-29 training run that
# used an adapted train_config with different patch/stride shapes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's A-cell first so the new date cell inherits the existing
# date-number-format style (keeps the same cellXf rather than minting a
# new one), then overwrite the value.
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 45167

$ws.Range("B4").Value = "train_config-230829-0.yml"
$ws.Range("E4").Value = "dataset03"
$ws.Range("F4").Value = "better performance metrics;"
$ws.Range("P4").Value = "cloud/pytorch-3dunet/resources/DW-3DUnet_lightsheet_boundary/named_copies/"
$ws.Range("Q4").Value = "cluster.s3it.uzh.ch:~/data/outputs/chpt-230829-0/"

# Move the active selection to the new row, matching the saved UI state.
$ws.Range("A4").Select() | Out-Null
